$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: free up the shared-string slots currently used by "NIS" and the
# long source description (rows 23/24) so the new strings we are about to
# introduce land *before* them in the shared-string table - matching the
# author's original authoring order (new rows were inserted above the old
# "NIS" / description rows, which then moved down to rows 29/30).
$ws.Range("A23").ClearContents()
$ws.Range("A24").ClearContents()

# Step 2: new header row (row 18) - bold/title style like row 9's header.
$ws.Range("B18").Value = "Number of employees"
$ws.Range("C18").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D18").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B18:D18").Style = "title"

# Step 3: Micro / Small / Medium / Large breakdown rows 19-22.
$ws.Range("A19").Value = "Micro"
$ws.Range("B19").Value = "<=10"
$ws.Range("C19").Value = "< USD50,000"
$ws.Range("D19").Value = " "

$ws.Range("A20").Value = "Small"
$ws.Range("B20").Value = "11-50"
$ws.Range("C20").Value = "USD50,000 - USD 250,000"
$ws.Range("D20").Value = " "

$ws.Range("A21").Value = "Medium"
$ws.Range("B21").Value = "51-100"
$ws.Range("C21").Value = "USD 250,000 - USD 500,000"
$ws.Range("D21").Value = " "

$ws.Range("A22").Value = "Large"
$ws.Range("B22").Value = ">100"
$ws.Range("C22").Value = ">USD 500,000"
$ws.Range("D22").Value = " "

$ws.Range("A19:D22").Style = "Normal"

# Step 4: restore the "NIS" / description rows, now at 29/30.
$ws.Range("A29").Value = "NIS"
$ws.Range("A29").Style = "title"

$ws.Range("A30").Value = "National Institute of Statistics (NIS), ""Nation-wide Establishment Listing of Cambodia 2009"", 2009, p. II-3, II-4. Available at http://catalog.ihsn.org/index.php/catalog/1496/related_materials"
$ws.Range("A30").Style = "source"
